# Update automatico via Actualizar 10-05-2020 19-57-13
# Appends 4 new daily rows (1-4 Oct 2020) to the "Condicion_Pacientes" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item("Condicion_Pacientes")

# New data rows: día, mes, año, Pruebas Realizadas, Pruebas Positivas
$newData = @(
    @(1, 10, 2020, 4638, 681),
    @(2, 10, 2020, 3969, 658),
    @(3, 10, 2020, 3351, 215),
    @(4, 10, 2020,  795, 215)
)

foreach ($rowData in $newData) {
    # Grow the table by one row - mirrors clicking "Add row" below the table.
    $lo.ListRows.Add() | Out-Null

    $r = $lo.Range.Row + $lo.Range.Rows.Count - 1

    # Copy the formatting of the row directly above so the new row keeps
    # the same cell styles (centered date, centered day/mes/año, ...).
    $ws.Range("A" + ($r - 1) + ":F" + ($r - 1)).Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)

    $ws.Cells.Item($r, 2).Value = $rowData[0]
    $ws.Cells.Item($r, 3).Value = $rowData[1]
    $ws.Cells.Item($r, 4).Value = $rowData[2]
    $ws.Cells.Item($r, 5).Value = $rowData[3]
    $ws.Cells.Item($r, 6).Value = $rowData[4]
    $ws.Cells.Item($r, 1).Formula = '=+Condicion_Pacientes[[#This Row],[día]]&"/"&Condicion_Pacientes[[#This Row],[mes]]&"/"&Condicion_Pacientes[[#This Row],[año]]'
}

$excel.CutCopyMode = $false

# Mirror the author's final viewport/selection state.
$lastRow = $lo.Range.Row + $lo.Range.Rows.Count - 1
$ws.Range("E" + ($lastRow + 2)).Select()
